$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the 'gender' column header to 'sex'
[void]$ws.Range("B1").Select()
$ws.Range("B1").Value = "sex"
